$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.611.03"
$ws.Range("E2").Value = '  +6.09%  '

$ws.Range("D3").Value = "'3.463.15"
$ws.Range("E3").Value = '  +4.52%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = "'415.41"
$ws.Range("E5").Value = '  +1.86%  '

$ws.Range("D6").Value = "'130.60"
$ws.Range("E6").Value = '  +18.05%  '

$ws.Range("D7").Value = "'3.456.08"
$ws.Range("E7").Value = '  +4.51%  '

$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = '  +1.54%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = "'0.690"
$ws.Range("E10").Value = '  +9.00%  '

$ws.Range("E11").Value = '  +29.94%  '

$ws.Range("D12").Value = "'43.84"
$ws.Range("E12").Value = '  +10.42%  '

$ws.Range("E13").Value = '  +0.45%  '

$ws.Range("D14").Value = "'4.028.22"
$ws.Range("E14").Value = '  +4.88%  '

$ws.Range("D15").Value = "'8.75"
$ws.Range("E15").Value = '  +3.97%  '

$ws.Range("D16").Value = "'20.32"
$ws.Range("E16").Value = '  +5.22%  '

$ws.Range("D17").Value = "'3.470.28"
$ws.Range("E17").Value = '  +4.75%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = "'62.577.84"
$ws.Range("E18").Value = '  +6.25%  '

$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").Value = "'1.05"
$ws.Range("E19").Value = '  +1.01%  '

$ws.Range("D20").Value = "'11.00"
$ws.Range("E20").Value = '  +2.74%  '

$ws.Range("D21").Value = "'0.0000137"
$ws.Range("E21").Value = '  +26.27%  '

$ws.Range("E22").Value = '  +1.87%  '

$ws.Range("D23").Value = "'13.24"
$ws.Range("E23").Value = '  +2.25%  '

$ws.Range("D24").Value = "'81.86"
$ws.Range("E24").Value = '  +9.37%  '

$ws.Range("D25").Value = "'313.34"
$ws.Range("E25").Value = '  +3.58%  '

$ws.Range("D26").Value = "'3.23"
$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").Value = "'30.58"
$ws.Range("E27").Value = '  +7.34%  '

$ws.Range("D28").Value = "'8.13"
$ws.Range("E28").Value = '  +3.72%  '

$ws.Range("D29").Value = "'7.79"
$ws.Range("E29").Value = '  +6.73%  '

$ws.Range("D30").Value = "'0.122"
$ws.Range("E30").Value = '  +8.72%  '

$ws.Range("D31").Value = "'0.178"
$ws.Range("E31").Value = '  +4.48%  '

$ws.Range("D32").Value = "'4.38"
$ws.Range("E32").Value = '  -1.85%  '

$ws.Range("D33").Value = "'44.93"
$ws.Range("E33").Value = '  +11.58%  '

$ws.Range("D34").Value = "'11.95"
$ws.Range("E34").Value = '  +4.87%  '

$ws.Range("D35").Value = "'2.63"
$ws.Range("E35").Value = '  +23.48%  '

$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").Value = "'0.0497"
$ws.Range("E37").Value = '  -5.91%  '

$ws.Range("D38").Value = "'52.67"
$ws.Range("E38").Value = '  +1.63%  '

$ws.Range("D39").Value = "'3.58"
$ws.Range("E39").Value = '  +3.31%  '

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = '  -6.57%  '

$ws.Range("D42").Value = "'2.01"
$ws.Range("E42").Value = '  +6.50%  '

$ws.Range("E43").Value = '  +3.17%  '

$ws.Range("D44").Value = "'137.38"
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").Value = "'17.75"
$ws.Range("E45").Value = '  +6.45%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = "'0.291"
$ws.Range("E46").Value = '  +5.02%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'4.00"
$ws.Range("E47").Value = '  +2.52%  '

$ws.Range("E48").Value = '  +0.56%  '

$ws.Range("D49").Value = "'22.61"
$ws.Range("E49").Value = '  +1.65%  '

$ws.Range("D50").Value = "'2.251.87"
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = "'3.820.08"
$ws.Range("E51").Value = '  +4.67%  '
